# Update countries & provincias Spain
# Refresh the daily COVID-19 stats table on sheet "Pais":
#  - update Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the countries whose figures moved
#  - a handful of neighbouring countries swapped ranking (and therefore row),
#    so their country name (column A) is corrected to match the new row order
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(4, 2).Value = 6581992
$ws.Cells.Item(4, 3).Value = 32517
$ws.Cells.Item(4, 4).Value = 3871100
$ws.Cells.Item(4, 5).Value = 2514823
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 831
$ws.Cells.Item(4, 8).Value = 196069

$ws.Cells.Item(11, 2).Value = 644438
$ws.Cells.Item(11, 3).Value = 2007
$ws.Cells.Item(11, 4).Value = 573003
$ws.Cells.Item(11, 5).Value = 56170
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 97
$ws.Cells.Item(11, 8).Value = 15265

$ws.Cells.Item(53, 2).Value = 58207
$ws.Cells.Item(53, 3).Value = 757
$ws.Cells.Item(53, 4).Value = 52284
$ws.Cells.Item(53, 5).Value = 5716
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 4
$ws.Cells.Item(53, 8).Value = 207

$ws.Cells.Item(68, 2).Value = 35603
$ws.Cells.Item(68, 3).Value = 143
$ws.Cells.Item(68, 4).Value = 22047
$ws.Cells.Item(68, 5).Value = 12944
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 5
$ws.Cells.Item(68, 8).Value = 612

$ws.Cells.Item(83, 2).Value = 18869
$ws.Cells.Item(83, 3).Value = 54
$ws.Cells.Item(83, 4).Value = 17884
$ws.Cells.Item(83, 5).Value = 866
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 119

$ws.Cells.Item(96, 2).Value = 9946
$ws.Cells.Item(96, 3).Value = 61
$ws.Cells.Item(96, 4).Value = 9100
$ws.Cells.Item(96, 5).Value = 783
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 63

$ws.Cells.Item(105, 2).Value = 7453
$ws.Cells.Item(105, 3).Value = 24
$ws.Cells.Item(105, 4).Value = 5635
$ws.Cells.Item(105, 5).Value = 1596
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 222

$ws.Cells.Item(110, 2).Value = 5655
$ws.Cells.Item(110, 3).Value = 2
$ws.Cells.Item(110, 4).Value = 3683
$ws.Cells.Item(110, 5).Value = 1796
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 176

$ws.Cells.Item(113, 1).Value = "Suazilandia"
$ws.Cells.Item(113, 2).Value = 4994
$ws.Cells.Item(113, 3).Value = 58
$ws.Cells.Item(113, 4).Value = 4103
$ws.Cells.Item(113, 5).Value = 793
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = 98

$ws.Cells.Item(114, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(114, 2).Value = 4990
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 4479
$ws.Cells.Item(114, 5).Value = 428
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 83

$ws.Cells.Item(119, 2).Value = 4747
$ws.Cells.Item(119, 3).Value = 11
$ws.Cells.Item(119, 4).Value = 1825
$ws.Cells.Item(119, 5).Value = 2860
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 62

$ws.Cells.Item(121, 1).Value = "Cuba"
$ws.Cells.Item(121, 2).Value = 4551
$ws.Cells.Item(121, 3).Value = 92
$ws.Cells.Item(121, 4).Value = 3779
$ws.Cells.Item(121, 5).Value = 666
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 2
$ws.Cells.Item(121, 8).Value = 106

$ws.Cells.Item(122, 1).Value = "Ruanda"
$ws.Cells.Item(122, 2).Value = 4479
$ws.Cells.Item(122, 3).Value = 19
$ws.Cells.Item(122, 4).Value = 2352
$ws.Cells.Item(122, 5).Value = 2105
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 22

$ws.Cells.Item(127, 1).Value = "Siria"
$ws.Cells.Item(127, 2).Value = 3416
$ws.Cells.Item(127, 3).Value = 65
$ws.Cells.Item(127, 4).Value = 797
$ws.Cells.Item(127, 5).Value = 2472
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 4
$ws.Cells.Item(127, 8).Value = 147

$ws.Cells.Item(128, 1).Value = "Eslovenia"
$ws.Cells.Item(128, 2).Value = 3389
$ws.Cells.Item(128, 3).Value = 77
$ws.Cells.Item(128, 4).Value = 2620
$ws.Cells.Item(128, 5).Value = 634
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 135

$ws.Cells.Item(129, 1).Value = "Mayotte"
$ws.Cells.Item(129, 2).Value = 3374
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 2964
$ws.Cells.Item(129, 5).Value = 370
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 40

$ws.Cells.Item(130, 1).Value = "Somalia"
$ws.Cells.Item(130, 2).Value = 3371
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 2738
$ws.Cells.Item(130, 5).Value = 536
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 97

$ws.Cells.Item(132, 1).Value = "Angola"
$ws.Cells.Item(132, 2).Value = 3217
$ws.Cells.Item(132, 3).Value = 125
$ws.Cells.Item(132, 4).Value = 1277
$ws.Cells.Item(132, 5).Value = 1810
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 4
$ws.Cells.Item(132, 8).Value = 130

$ws.Cells.Item(133, 1).Value = "Lituania"
$ws.Cells.Item(133, 2).Value = 3199
$ws.Cells.Item(133, 3).Value = 36
$ws.Cells.Item(133, 4).Value = 2030
$ws.Cells.Item(133, 5).Value = 1083
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 86

$ws.Cells.Item(134, 1).Value = "Sri Lanka"
$ws.Cells.Item(134, 2).Value = 3155
$ws.Cells.Item(134, 3).Value = 8
$ws.Cells.Item(134, 4).Value = 2955
$ws.Cells.Item(134, 5).Value = 188
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 12

$ws.Cells.Item(145, 1).Value = "Botsuana"
$ws.Cells.Item(145, 2).Value = 2252
$ws.Cells.Item(145, 3).Value = 126
$ws.Cells.Item(145, 4).Value = 546
$ws.Cells.Item(145, 5).Value = 1696
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 10

$ws.Cells.Item(146, 1).Value = "Benin"
$ws.Cells.Item(146, 2).Value = 2242
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 1793
$ws.Cells.Item(146, 5).Value = 409
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 40

$ws.Cells.Item(147, 1).Value = "Malta"
$ws.Cells.Item(147, 2).Value = 2204
$ws.Cells.Item(147, 3).Value = 42
$ws.Cells.Item(147, 4).Value = 1803
$ws.Cells.Item(147, 5).Value = 387
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 14

$ws.Cells.Item(148, 1).Value = "Islandia"
$ws.Cells.Item(148, 2).Value = 2157
$ws.Cells.Item(148, 3).Value = 4
$ws.Cells.Item(148, 4).Value = 2072
$ws.Cells.Item(148, 5).Value = 75
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 10

$ws.Cells.Item(149, 1).Value = "Birmania"
$ws.Cells.Item(149, 2).Value = 2150
$ws.Cells.Item(149, 3).Value = 261
$ws.Cells.Item(149, 4).Value = 625
$ws.Cells.Item(149, 5).Value = 1511
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 2
$ws.Cells.Item(149, 8).Value = 14

$ws.Cells.Item(154, 1).Value = "Guyana"
$ws.Cells.Item(154, 2).Value = 1750
$ws.Cells.Item(154, 3).Value = 47
$ws.Cells.Item(154, 4).Value = 1088
$ws.Cells.Item(154, 5).Value = 613
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 1
$ws.Cells.Item(154, 8).Value = 49

$ws.Cells.Item(155, 1).Value = "Uruguay"
$ws.Cells.Item(155, 2).Value = 1741
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 1478
$ws.Cells.Item(155, 5).Value = 218
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 45

$ws.Cells.Item(168, 2).Value = 901
$ws.Cells.Item(168, 3).Value = 3
$ws.Cells.Item(168, 4).Value = 866
$ws.Cells.Item(168, 5).Value = 20
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 15

$ws.Cells.Item(177, 2).Value = 469
$ws.Cells.Item(177, 3).Value = 3
$ws.Cells.Item(177, 4).Value = 374
$ws.Cells.Item(177, 5).Value = 94
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 1

$ws.Cells.Item(188, 2).Value = 208
$ws.Cells.Item(188, 3).Value = 1
$ws.Cells.Item(188, 4).Value = 204
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 1

$ws.Cells.Item(189, 2).Value = 180
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 157
$ws.Cells.Item(189, 5).Value = 16
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 7

$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0
